$wb = $excel.ActiveWorkbook
$wsGlobal = $wb.Worksheets.Item(1)

# Add a new worksheet "Aciclovir" right after "Global", containing what used to be
# the Aciclovir parameter row (previously row 2 of "Global").
$wsAciclovir = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsGlobal)
$wsAciclovir.Name = "Aciclovir"

$wsAciclovir.Range("A1").Value = "Container Path"
$wsAciclovir.Range("B1").Value = "Parameter Name"
$wsAciclovir.Range("C1").Value = "Value"
$wsAciclovir.Range("D1").Value = "Units"

$wsAciclovir.Range("A2").Value = "Aciclovir"
$wsAciclovir.Range("B2").Value = "Lipophilicity"
$wsAciclovir.Range("C2").Value = -0.1
$wsAciclovir.Range("D2").Value = "Log Units"
$wsAciclovir.Range("C1:C2").NumberFormat = "0.0000"

[void]$wsAciclovir.Rows("1:2").Select()

# Update the "Global" sheet: row 2 now holds the "EHC continuous fraction" /
# "Organism|Liver" parameter instead of the Aciclovir/Lipophilicity one, and the
# Units column (D2) is no longer used.
$wsGlobal.Range("B2").Value = "EHC continuous fraction"
$wsGlobal.Range("A2").Value = "Organism|Liver"
$wsGlobal.Range("C2").Value = 1
$wsGlobal.Range("D2").Value = ""

[void]$wsGlobal.Activate()
[void]$wsGlobal.Range("C10").Select()
